# chore: adapt column header formatting to respective input file names
#
# Renames the "_old"/"_new" header-name suffixes to the respective
# format-version suffixes ("_FV2310" / "_FV2404"), wraps the sheet's data
# range in an Excel Table ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row (row 1) cells: "<name>_old" -> "<name>_FV2310"
#    and "<name>_new" -> "<name>_FV2404". Column K ("diff") stays as-is.
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $current = $cell.Value()
    $cell.Value = $current.Replace("_old", "_FV2310")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $current = $cell.Value()
    $cell.Value = $current.Replace("_new", "_FV2404")
}

# 2) Turn the used range into an Excel Table ("Table1") so the header row
#    doubles as filter buttons, using the renamed headers as column names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U73"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split beneath row 1).
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$ws.Range("A1").Select()
